$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 15:52"

# Row 9 - Valencia/Valencia
$ws.Range("C9").Value = 1717
$ws.Range("D9").Value = 2668
$ws.Range("E9").Value = 464

# Row 13 - Alacant/Alicante
$ws.Range("C13").Value = 1314
$ws.Range("D13").Value = 1707
$ws.Range("E13").Value = 357

# Row 33 - Tenerife
$ws.Range("B33").Value = 1249
$ws.Range("C33").Value = 370
$ws.Range("D33").Value = 807

# Row 36 - Castello/Castellon
$ws.Range("C36").Value = 329
$ws.Range("D36").Value = 739
$ws.Range("E36").Value = 124

# Row 50 - Gran Canaria
$ws.Range("C50").Value = 196
$ws.Range("D50").Value = 233

# Row 56 - La Palma
$ws.Range("C56").Value = 20
$ws.Range("D56").Value = 59

# Row 57 - Lanzarote
$ws.Range("B57").Value = 77
$ws.Range("D57").Value = 55

# Row 59 - Fuerteventura
$ws.Range("C59").Value = 26
$ws.Range("D59").Value = 16
